# Conserto do erro com o rótulo da coluna 2050 nas tabelas e
# retirada das linhas com total das tabelas.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# Helper: write a text label into a cell while preserving the header
# cell's existing style (bold / border / centered) instead of letting
# Excel re-interpret a numeric-looking string as a number.
function Set-TextLabel($ws, $cellAddr, $styleSourceAddr, $text) {
    $ws.Range($cellAddr).Value = "'" + $text
    $ws.Range($styleSourceAddr).Copy()
    $ws.Range($cellAddr).PasteSpecial($xlPasteFormats)
}

# Sheets 1-3: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)" - column E header should read "2050", and the
# bogus "Total" row (row 13) is removed.
foreach ($idx in 1..3) {
    $ws = $wb.Worksheets.Item($idx)
    Set-TextLabel $ws "E1" "D1" "2050"
    $ws.Rows.Item(13).Delete()
}

# Sheet 4: "Potencia Incremental - SIN(MW)" uses year-range labels, so the
# column E header becomes "2041-2050" instead of plain "2050".
$ws4 = $wb.Worksheets.Item(4)
Set-TextLabel $ws4 "E1" "D1" "2041-2050"
$ws4.Rows.Item(13).Delete()

# Sheet 5: "Emissoes Totais (MtCO2eq)" only needs the column E header fixed;
# it never had a "Total" row to remove.
$ws5 = $wb.Worksheets.Item(5)
Set-TextLabel $ws5 "E1" "D1" "2050"

# Sheet 6: "Custo Total (bilhões de R$)" has no year columns, just drop the
# trailing "Total" row (row 4).
$ws6 = $wb.Worksheets.Item(6)
$ws6.Rows.Item(4).Delete()
